$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-357) holds a "Förändrad" (changed) date that was bumped
# from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188) for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 357 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
